# Apply "Add smoke test case for update revision" edit to the FIRevisions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FIRevisions")

# --- New annotations in column E (a "ParentLocator" style notes column) ---
$ws.Range("E1").Value = "ParentLocator"
$ws.Range("E9").Value  = '//*[@id=''Contract_listbox'']'
$ws.Range("E10").Value = '//*[@id=''Category_listbox'']'
$ws.Range("E13").Value = '//*[@id="cmbChartOfAccount_listbox"]'
$ws.Range("E14").Value = '//*[@id="ddlLeaseTypes_listbox"]'
$ws.Range("E15").Value = '//*[@id="ddlChartOfAccountTypes_listbox"]'
$ws.Range("E16").Value = '//*[@id="ddlMaintainCodeTypes_listbox"]'

# E4 carries the same (empty-value) formatting as the locator-style cells above it.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Normalize xpath locator quoting style (double quotes -> single quotes) ---
$ws.Range("B33").Value = '//*[@id=''Prepayment'']'
$ws.Range("B34").Value = '//*[@id=''ROU_InitialDirectCosts'']'
$ws.Range("B35").Value = '//*[@id=''ROU_LandlordAllowance'']'

# --- New smoke test rows for the FASB classification dropdown + notes field ---
$ws.Range("A59").Value = "FASBClassificationType"
$ws.Range("B59").Value = '//*[@aria-owns="FASBClassificationTypeID_listbox"]'
$ws.Range("C59").Value = "by_xpath"
$ws.Range("E59").Value = '//*[@id="FASBClassificationTypeID_listbox"]'

$ws.Range("A60").Value = "FASBClassificationTypeParent"
$ws.Range("B60").Value = '//*[@id="FASBClassificationTypeID_listbox"]'
$ws.Range("C60").Value = "by_xpath"

$ws.Range("A61").Value = "txtNotes"
$ws.Range("B61").Value = "Notes_Description"
$ws.Range("C61").Value = "by_id"

# B61 picks up the same highlighted/monospace style used by other locator cells.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B61").Value = "Notes_Description"

# --- Column width tweaks (B/D narrower, new column E) ---
$ws.Columns.Item(2).ColumnWidth = 56.92
$ws.Columns.Item(4).ColumnWidth = 45.25
$ws.Columns.Item(5).ColumnWidth = 41.75

# --- Update the view so the frozen pane / selection matches the new row count ---
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A26").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B41").Select() | Out-Null
